# Add 7 new rows (205-211) of landscaping observations for 2025-06-08
# to the bottom of Sheet1, continuing the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the date's number format (m/d/yyyy, same style as the rows above)
# down into the new rows by copying the format from the last existing
# data row (204) instead of inventing a new style.
$ws.Range("A204").Copy()
$ws.Range("A205:A211").PasteSpecial(-4122)

# New rows data. Columns:
# A Date(serial) | B Plant_Type | C Plant_Size | D Low | E High | F Temp_Diff(formula)
# G Rain | H Growth | I Pruned | J Quadrant | K Shade | L UV | M Humidity
# N Dew_Point | O Pressure | P Wind_Gust | Q Cloud_Cover | R Visibility | S AQI | T Pollen
$rows = @(
    @{R=205; A=45816; B="Flowering";    C="Large";  D=63; E=68; G=1.31; H=0.35; I="No"; J=2; K="Dark";    L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29},
    @{R=206; A=45816; B="Nonflowering"; C="Medium"; D=63; E=68; G=1.31; H=0.5;  I="No"; J=3; K="Dark";    L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29},
    @{R=207; A=45816; B="Nonflowering"; C="Small";  D=63; E=68; G=1.31; H=0.75; I="No"; J=3; K="Neutral"; L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29},
    @{R=208; A=45816; B="Nonflowering"; C="Medium"; D=63; E=68; G=1.31; H=1.5;  I="No"; J=3; K="Neutral"; L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29},
    @{R=209; A=45816; B="Nonflowering"; C="Medium"; D=63; E=68; G=1.31; H=2.35; I="No"; J=3; K="Bright";  L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29},
    @{R=210; A=45816; B="Nonflowering"; C="Large";  D=63; E=68; G=1.31; H=0.75; I="No"; J=4; K="Dark";    L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29},
    @{R=211; A=45816; B="Tree";         C="Medium"; D=63; E=68; G=1.31; H=3.75; I="No"; J=1; K="Neutral"; L=7; M=1; N=68; O=29.64; P=13; Q=0.91; R_=3; S=100; T=29}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R_
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}

# Column F (Temp_Diff) continues the existing ABS(D-E) shared formula pattern.
$ws.Range("F205:F211").Formula = "=ABS(D205-E205)"

# Move the view/selection the way the source file ended up (best effort;
# some window-chrome state isn't persisted by this host).
$ws.Range("Z2").Select() | Out-Null
